$d = $word.ActiveDocument

function FindRange([string]$text) {
    $r = $d.Content
    $r.Find.ClearFormatting()
    $ok = $r.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Could not find text: $text"
    }
    return $r
}

function ReplaceText([string]$old, [string]$new) {
    $r = $d.Content
    $r.Find.ClearFormatting()
    $ok = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        throw "Could not find/replace text: $old"
    }
}

# --- Change 1: Dedication paragraph - drop bold formatting on children/nephew/nieces
#     names and merge all those runs into the surrounding italic-only text. ---
$beforeNames = FindRange "For my children"
$afterFirstRunPos = $beforeNames.End

$beforeWish = FindRange ". My wish is"
$wishStart = $beforeWish.Start

$namesRange = $d.Range($afterFirstRunPos, $wishStart)
if ($namesRange.Text -ne " Rosie, Joey, and Zach; my nephew Elliott; and my nieces Amy and Lyla") {
    throw "Unexpected dedication text: $($namesRange.Text)"
}
$namesRange.Delete()

$insPoint = $d.Range($afterFirstRunPos, $afterFirstRunPos)
$insPoint.InsertAfter(" Rosie, Joey, and Zach; my nephew Elliott; and my nieces Amy and Lyla")

# --- Change 2: Acknowledgements opening paragraph reworded ---
ReplaceText `
    "No-one really knows how hard a PhD is until they are already well beyond the point of no return. It is demanding, challenging and often thankless and lonely work, working countless hours, striving to find meaning among masses of data or iterating endlessly to concisely express complex, nebulous and elusive ideas. It is an endeavour made even harder when one is a mature student with financial and parental responsibilities. At times, especially during the final unfunded writing-up period, the impacts upon my life and those around me have been huge and unreasonable. For this reason, the greatest thanks of all go to my wife" `
    "No-one really knows how hard a PhD is until they are already well beyond the point of no return. It is demanding, challenging and often thankless and lonely work. You spend countless hours striving to find meaning among masses of data or iterating endlessly to concisely express complex, nebulous and elusive ideas. This endeavour is made even harder when you are a mature student with financial and parental responsibilities. At times, especially during the final unfunded writing-up period, the impacts upon my life and those around me have been huge and unreasonable. For this reason, the greatest thanks of all go to my wife"

# --- Change 3: Rob Comber bullet - add Matt Wood ---
ReplaceText `
    "Rob Comber, Simon Bowen and all the other lecturers who taught me valuable Digital Civics and research skills during the MRes, that helped shape me into the researcher I am." `
    "Rob Comber, Simon Bowen, Matt Wood and all the other lecturers who taught me valuable Digital Civics and research skills during the MRes, that helped shape me into the researcher I am."

# --- Change 4: Stuart Wheater bullet reworded ---
ReplaceText `
    "Stuart Wheater, for myriad data discussions and tactical discussions through a difficult period that often went beyond project business." `
    "Stuart Wheater, for myriad data discussions and tactical discussions that helped me get through a difficult period and often went beyond project business."

# --- Change 5: Marshall bullet - move ", and supporting me in adapting it" to after the citation ---
$genRange = FindRange "document generation"
$genEnd = $genRange.End

$fullPhrase = FindRange "document generation, and supporting me in adapting it"
$phraseEnd = $fullPhrase.End

$toDelete = $d.Range($genEnd, $phraseEnd)
if ($toDelete.Text -ne ", and supporting me in adapting it") {
    throw "Unexpected Marshall text: $($toDelete.Text)"
}
$toDelete.Delete()

$beforeThisMade = FindRange "; this made thesis development"
$insPos = $beforeThisMade.Start
$insRange2 = $d.Range($insPos, $insPos)
$insRange2.InsertBefore(", and supporting me in adapting it")

Write-Host "All changes applied"
